$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Get-ParaIndexContaining([string]$needle) {
    $probe = $d.Content.Duplicate
    $null = $probe.Find.Execute($needle)
    return $probe.Paragraphs(1).Index
}

# ------------------------------------------------------------------
# 1. Remove the _GoBack bookmark from the end of the
#    "Print the bracket ... using PLA" paragraph (it gets re-added
#    later, at the end of the new final paragraph).
# ------------------------------------------------------------------
$plaIdx = Get-ParaIndexContaining("Print the bracket at 200um spacing")
$plaRange = $d.Paragraphs($plaIdx).Range
$plaXml = '<w:p ' + $wns + '>' + `
    '<w:r><w:t>Print the bracket at 200um spacing</w:t></w:r>' + `
    '<w:r><w:t>, using PLA</w:t></w:r>' + `
    '</w:p>'
$plaRange.InsertXML($plaXml)

# ------------------------------------------------------------------
# 2. Insert a new bulleted "magnet" list item before the
#    "At step 4, fit the magnet ..." paragraph, matching the
#    existing bullet list ("4mm rod ..." / "A 4 to 5mm ...").
# ------------------------------------------------------------------
$stepIdx = Get-ParaIndexContaining("At step 4, fit the magnet")
$d.Paragraphs($stepIdx).Range.InsertParagraphBefore()

$newListPara = $d.Paragraphs($stepIdx)
$newListPara.Range.Text = "magnet"
$newListPara.Style = "List Paragraph"

$templateListPara = $d.Paragraphs($stepIdx - 1)
$newListPara.Range.ListFormat.ListTemplate = $templateListPara.Range.ListFormat.ListTemplate
$newListPara.Range.ListFormat.ListLevelNumber = $templateListPara.Range.ListFormat.ListLevelNumber

# ------------------------------------------------------------------
# 3. Rewrite the "At step 4 ..." paragraph: "expoxy" -> "hot glue"
#    (dropping the spell-check proofErr markers around it) and append
#    the new sentence about testing the hall effect switch.
# ------------------------------------------------------------------
$stepIdx2 = Get-ParaIndexContaining("At step 4, fit the magnet")
$stepRange = $d.Paragraphs($stepIdx2).Range
$stepXml = '<w:p ' + $wns + '>' + `
    '<w:r><w:t xml:space="preserve">At step 4, fit the magnet to the bottom of the deck using </w:t></w:r>' + `
    '<w:r><w:t>hot glue</w:t></w:r>' + `
    '<w:r><w:t>, ensure that is position in the centre with North South poles in line with the deck (not across the deck)</w:t></w:r>' + `
    '<w:r><w:t>. Before gluing test with the hall effect switch, it is only activated by one pole of the magnet!</w:t></w:r>' + `
    '</w:p>'
$stepRange.InsertXML($stepXml)

# ------------------------------------------------------------------
# 4. Add a new paragraph after it with the home-switch guidance, and
#    re-home the _GoBack bookmark at its end.
# ------------------------------------------------------------------
$stepIdx3 = Get-ParaIndexContaining("At step 4, fit the magnet")
$d.Paragraphs($stepIdx3).Range.InsertParagraphAfter()

$homeRange = $d.Paragraphs($stepIdx3 + 1).Range
$homeXml = '<w:p ' + $wns + '>' + `
    '<w:r><w:t>Choose a position for the home switch so that the power pickups are in the middle of their travel.</w:t></w:r>' + `
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' + `
    '</w:p>'
$homeRange.InsertXML($homeXml)

$d.Save()
